# Update the "birimler" (units) workbook:
#  - MEZITLI sheet gets a brand-new, shorter village list (A1:A23),
#    the old rows A24:A73 are cleared out.
#  - MUHENDIS sheet gains a new engineer "Volkan IZCI" at the end.
#  - TEKNIKER sheet gains a new technician "Ayse PEKCETIN" inserted
#    as the 5th row (existing rows shift down).
#  - Selections on a few sheets are updated to match where the author
#    was last working.

$wb = $excel.ActiveWorkbook

# ---- TEKNIKER (6th tab) ------------------------------------------------
# New technician "Ayse PEKCETIN" is inserted as the 5th row; the rest
# of the list shifts down to make room.
$tekniker = $wb.Worksheets.Item(6)
$tekniker.Rows.Item(5).Insert()
$tekniker.Cells.Item(5, 1).Value = "Ayşe PEKÇETİN"

# ---- MEZITLI (2nd tab) -----------------------------------------------
$mezitli = $wb.Worksheets.Item(2)

$mezitliVillages = @(
    "Akarca",
    "Bozön",
    "Çamlıca",
    "Çankaya",
    "Çevlik",
    "Davultepe",
    "Demirışık",
    "Doğançay",
    "Doğlu",
    "Fındıkpınarı",
    "Kale",
    "Kaleburnu",
    "Kocayer",
    "Kuyuluk",
    "Kuzucu",
    "Kuzucubelen",
    "Mezitli",
    "Pelitkoyağı",
    "Sarılar",
    "Takanlı",
    "Tece",
    "Tepeköy",
    "Tolköy"
)

for ($i = 0; $i -lt $mezitliVillages.Length; $i++) {
    $mezitli.Cells.Item($i + 1, 1).Value = $mezitliVillages[$i]
}

# The sheet used to hold 73 villages; only 23 remain, so wipe the rest.
$mezitli.Range("A24:A73").ClearContents()

# Shrink the column back down now that the longest entries are gone.
$mezitli.Columns.Item(1).AutoFit()

# ---- MUHENDIS (5th tab) ------------------------------------------------
$muhendis = $wb.Worksheets.Item(5)
$muhendis.Cells.Item(11, 1).Value = "Volkan İZCİ"

# ---- Update the last-used cell on each sheet to match where the ------
# ---- author ended up, finishing on TEKNIKER so it stays the active ---
# ---- tab (as in the original file).                                 --
$mezitli.Range("E10").Select()
$wb.Worksheets.Item(4).Range("I42").Select()
$muhendis.Range("A12").Select()
$tekniker.Range("B17").Select()
